$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# Small single-line number tweaks in the stack trace
Replace-Text "JavaMethodService.java:163" "JavaMethodService.java:162"
Replace-Text "AbstractService.java:136" "AbstractService.java:135"
Replace-Text "EvaluationServices.call(EvaluationServices.java:168)" "EvaluationServices.call(EvaluationServices.java:172)"
Replace-Text "AstEvaluator.caseCall(AstEvaluator.java:189)" "AstEvaluator.caseCall(AstEvaluator.java:186)"
Replace-Text "AstEvaluator.caseCall(AstEvaluator.java:183)" "AstEvaluator.caseCall(AstEvaluator.java:180)"
Replace-Text "AstEvaluator.eval(AstEvaluator.java:112)" "AstEvaluator.eval(AstEvaluator.java:109)"
Replace-Text "GeneratedMethodAccessor74.invoke" "GeneratedMethodAccessor73.invoke"

# Two occurrences of AstSwitch.java:118 -> 119, replace one at a time (search forward each time)
Replace-Text "AstSwitch.doSwitch(AstSwitch.java:118)" "AstSwitch.doSwitch(AstSwitch.java:119)"
Replace-Text "AstSwitch.doSwitch(AstSwitch.java:118)" "AstSwitch.doSwitch(AstSwitch.java:119)"

# Replace the tail block of the stack trace (maven/surefire/equinox launcher lines)
# with the jdt junit4 runner lines.
$oldBlock = "`tat org.apache.maven.surefire.junit4.JUnit4Provider.execute(JUnit4Provider.java:264)`n" +
            "`tat org.apache.maven.surefire.junit4.JUnit4Provider.executeTestSet(JUnit4Provider.java:153)`n" +
            "`tat org.apache.maven.surefire.junit4.JUnit4Provider.invoke(JUnit4Provider.java:124)`n" +
            "`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n" +
            "`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n" +
            "`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n" +
            "`tat java.lang.reflect.Method.invoke(Method.java:498)`n" +
            "`tat org.apache.maven.surefire.util.ReflectionUtils.invokeMethodWithArray2(ReflectionUtils.java:208)`n" +
            "`tat org.apache.maven.surefire.booter.ProviderFactory`$ProviderProxy.invoke(ProviderFactory.java:156)`n" +
            "`tat org.apache.maven.surefire.booter.ProviderFactory.invokeProvider(ProviderFactory.java:82)`n" +
            "`tat org.eclipse.tycho.surefire.osgibooter.OsgiSurefireBooter.run(OsgiSurefireBooter.java:91)`n" +
            "`tat org.eclipse.tycho.surefire.osgibooter.HeadlessTestApplication.run(HeadlessTestApplication.java:21)`n" +
            "`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n" +
            "`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n" +
            "`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n" +
            "`tat java.lang.reflect.Method.invoke(Method.java:498)`n" +
            "`tat org.eclipse.equinox.internal.app.EclipseAppContainer.callMethodWithException(EclipseAppContainer.java:587)`n" +
            "`tat org.eclipse.equinox.internal.app.EclipseAppHandle.run(EclipseAppHandle.java:198)`n" +
            "`tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.runApplication(EclipseAppLauncher.java:134)`n" +
            "`tat org.eclipse.core.runtime.internal.adaptor.EclipseAppLauncher.start(EclipseAppLauncher.java:104)`n" +
            "`tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:388)`n" +
            "`tat org.eclipse.core.runtime.adaptor.EclipseStarter.run(EclipseStarter.java:243)`n" +
            "`tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)`n" +
            "`tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)`n" +
            "`tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)`n" +
            "`tat java.lang.reflect.Method.invoke(Method.java:498)`n" +
            "`tat org.eclipse.equinox.launcher.Main.invokeFramework(Main.java:656)`n" +
            "`tat org.eclipse.equinox.launcher.Main.basicRun(Main.java:592)`n" +
            "`tat org.eclipse.equinox.launcher.Main.run(Main.java:1498)`n" +
            "`tat org.eclipse.equinox.launcher.Main.main(Main.java:1471)"

$newBlock = "`tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)`n" +
            "`tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)`n" +
            "`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:538)`n" +
            "`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:760)`n" +
            "`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:460)`n" +
            "`tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:206)"

Replace-Text $oldBlock $newBlock
